$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New player roster table (rows 2-19), columns A=Player, B=Position, C=Team
$data = @(
    @("Chris Paul",        "PG",       "San Antonio Spurs"),
    @("Dillon Brooks",     "SG,SF",    "Houston Rockets"),
    @("Payton Pritchard",  "PG,SG",    "Boston Celtics"),
    @("Jalen Green",       "PG,SG",    "Houston Rockets"),
    @("Paolo Banchero",    "SF,PF",    "Orlando Magic"),
    @("Deni Avdija",       "SF,PF",    "Portland Trail Blazers"),
    @("Ayo Dosunmu",       "PG,SG,SF", "Chicago Bulls"),
    @("Pascal Siakam",     "SF,PF,C",  "Indiana Pacers"),
    @("Nikola Jokic",      "C",        "Denver Nuggets"),
    @("Rudy Gobert",       "C",        "Minnesota Timberwolves"),
    @("Jakob Poeltl",      "C",        "Toronto Raptors"),
    @("Stephon Castle",    "PG,SG",    "San Antonio Spurs"),
    @("Jaylen Brown",      "SG,SF",    "Boston Celtics"),
    @("Ty Jerome",         "PG,SG",    "Cleveland Cavaliers"),
    @("Khris Middleton",   "SF",       "Milwaukee Bucks"),
    @("Chet Holmgren",     "PF,C",     "Oklahoma City Thunder"),
    @("Jalen Suggs",       "PG,SG",    "Orlando Magic"),
    @("Russell Westbrook", "PG,SG",    "Denver Nuggets")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row++
}
